$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# like "1.41" or "0.550" keep their exact textual representation
# instead of being coerced into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "93.552.89"
$ws.Range("E2").Value = "  -4.19%  "
$ws.Range("D3").Value = "3.409.36"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "234.93"
$ws.Range("E5").Value = "  -6.93%  "
$ws.Range("D6").Value = "634.67"
$ws.Range("E6").Value = "  -3.83%  "
$ws.Range("D7").Value = "1.41"
$ws.Range("E7").Value = "  -1.65%  "
$ws.Range("D8").Value = "0.389"
$ws.Range("E8").Value = "  -8.59%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "0.950"
$ws.Range("E10").Value = "  -6.87%  "
$ws.Range("D11").Value = "3.399.97"
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("E12").Value = "  -5.70%  "
$ws.Range("D13").Value = "41.42"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").Value = "6.05"
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").Value = "93.103.71"
$ws.Range("E15").Value = "  -4.44%  "
$ws.Range("D16").Value = "4.040.23"
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("D17").Value = "0.0000248"
$ws.Range("E17").Value = "  -2.55%  "
$ws.Range("D18").Value = "8.27"
$ws.Range("E18").Value = "  -6.90%  "
$ws.Range("D19").Value = "3.388.21"
$ws.Range("E19").Value = "  +0.68%  "
$ws.Range("D20").Value = "17.41"
$ws.Range("E20").Value = "  -2.48%  "
$ws.Range("D21").Value = "11.15"
$ws.Range("E21").Value = "  +3.75%  "
$ws.Range("D22").Value = "0.483"
$ws.Range("E22").Value = "  -13.09%  "
$ws.Range("D23").Value = "491.78"
$ws.Range("E23").Value = "  -3.98%  "
$ws.Range("D24").Value = "3.16"
$ws.Range("E24").Value = "  -5.54%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "0.0000187"
$ws.Range("E25").Value = "  -6.96%  "
$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").Value = "6.36"
$ws.Range("E26").Value = "  -4.10%  "
$ws.Range("D27").Value = "90.28"
$ws.Range("E27").Value = "  -7.48%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "11.88"
$ws.Range("E28").Value = "  -3.19%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "3.562.31"
$ws.Range("E29").Value = "  +0.97%  "
$ws.Range("D30").Value = "11.55"
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("D32").Value = "2.69"
$ws.Range("E32").Value = "  +5.10%  "
$ws.Range("D33").Value = "0.134"
$ws.Range("E33").Value = "  -9.04%  "
$ws.Range("D34").Value = "0.178"
$ws.Range("E34").Value = "  -6.99%  "
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").Value = "29.78"
$ws.Range("E36").Value = "  +3.92%  "
$ws.Range("D37").Value = "0.550"
$ws.Range("E37").Value = "  -2.70%  "
$ws.Range("D38").Value = "540.51"
$ws.Range("E38").Value = "  +4.71%  "
$ws.Range("D39").Value = "7.49"
$ws.Range("E39").Value = "  -5.68%  "
$ws.Range("D40").Value = "1.43"
$ws.Range("E40").Value = "  -5.15%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "0.915"
$ws.Range("E42").Value = "  +7.94%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "0.148"
$ws.Range("E43").Value = "  -2.94%  "
$ws.Range("D44").Value = "24.03"
$ws.Range("E44").Value = "  -1.60%  "
$ws.Range("D45").Value = "1.66"
$ws.Range("E45").Value = "  -2.98%  "
$ws.Range("D46").Value = "0.0407"
$ws.Range("E46").Value = "  -6.27%  "
$ws.Range("D47").Value = "5.47"
$ws.Range("E47").Value = "  -4.12%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "2.11"
$ws.Range("E48").Value = "  +4.88%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "52.63"
$ws.Range("E49").Value = "  -4.06%  "
$ws.Range("B50").Value = "MantraDAO"
$ws.Range("C50").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D50").Value = "3.41"
$ws.Range("E50").Value = "  -6.42%  "
$ws.Range("D51").Value = "3.13"
$ws.Range("E51").Value = "  -0.66%  "
